$d = $word.ActiveDocument

# Locate the paragraph holding the "_____________ <tab> Хитрин Артём Сергеевич"
# signature line by its distinctive surname text (robust to absolute offsets).
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Хитрин*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target signature paragraph containing 'Хитрин'"
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$pStart = $targetPara.Range.Start
$pEnd = $targetPara.Range.End

# Remove all of the paragraph's run content (the underscores, tab and name),
# leaving only the paragraph mark itself.
$textRange = $d.Range($pStart, $pEnd - 1)
if ($textRange.Start -lt $textRange.End) {
    $textRange.Delete()
}

# Mark the now-empty paragraph's mark as English (US).
$emptiedPara = $d.Paragraphs.Item($targetIndex)
$emptiedPara.Range.LanguageID = "en-US"

# The following (already-empty) paragraph gets the same paragraph-mark language.
$nextPara = $d.Paragraphs.Item($targetIndex + 1)
$nextPara.Range.LanguageID = "en-US"
